# "implemented pad list with connected nets"
# Row 3 (B3:I3) held the pad numbers 8..1 counting down across the pad
# list; renumber it to 7..0 (each value decremented by one) and leave the
# active selection on A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0

$ws.Range("A3").Select()
